$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: update title (D12) and link (E12)
$ws.Range("D12").Value = "“케라스 창시자에게 배우는 딥러닝 2판”이 출간되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/08/16/%ec%bc%80%eb%9d%bc%ec%8a%a4-%ec%b0%bd%ec%8b%9c%ec%9e%90%ec%97%90%ea%b2%8c-%eb%b0%b0%ec%9a%b0%eb%8a%94-%eb%94%a5%eb%9f%ac%eb%8b%9d-2%ed%8c%90%ec%9d%b4-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88%ec%8a%b5/"

# Row 27: update title (D27) and link (E27)
$ws.Range("D27").Value = "AWS Inferentia를 이용한 모델 서빙 비용 최적화: 모델 서버 비용 2배 줄이기 2탄"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/aws-inferentia-2/"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "[pandas] SQL의 JOIN 역할을 하는 pd.merge와 df.join"
$ws.Range("E51").Value = "https://bskyvision.com/entry/pandas-SQL%EC%9D%98-JOIN-%EC%97%AD%ED%95%A0%EC%9D%84-%ED%95%98%EB%8A%94-pdmerge%EC%99%80-dfjoin"
